$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued cells whose new content looks numeric: force text storage
# by setting NumberFormat to "@" before assigning the value, then
# resetting the style back to Normal so no stray formatting is left behind.
function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "61.758.17"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "3.402.69"
$ws.Range("E3").Value = "  +1.12%  "
Set-TextValue $ws "D4" "0.998"
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue $ws "D5" "579.02"
$ws.Range("E5").Value = "  +1.32%  "
Set-TextValue $ws "D6" "137.86"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.403.67"
$ws.Range("E8").Value = "  +1.19%  "
Set-TextValue $ws "D9" "0.474"
$ws.Range("E9").Value = "  -0.12%  "
Set-TextValue $ws "D10" "7.54"
$ws.Range("E10").Value = "  -0.75%  "
Set-TextValue $ws "D11" "0.126"
$ws.Range("E11").Value = "  +2.95%  "
Set-TextValue $ws "D12" "0.390"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "3.969.21"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("E14").Value = "  +2.25%  "
Set-TextValue $ws "D15" "0.0000177"
$ws.Range("E15").Value = "  +2.95%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws "D16" "25.97"
$ws.Range("E16").Value = "  +3.17%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.396.33"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "61.759.73"
$ws.Range("E18").Value = "  +1.24%  "
Set-TextValue $ws "D19" "14.25"
$ws.Range("E19").Value = "  +2.46%  "
Set-TextValue $ws "D20" "5.89"
$ws.Range("E20").Value = "  +1.39%  "
Set-TextValue $ws "D21" "9.48"
$ws.Range("E21").Value = "  +0.51%  "
Set-TextValue $ws "D22" "377.86"
$ws.Range("E22").Value = "  +1.45%  "
Set-TextValue $ws "D23" "0.560"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").Value = "3.527.49"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws "D25" "0.0000127"
$ws.Range("E25").Value = "  +8.85%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D26" "1.00"
$ws.Range("E26").Value = "  +0.02%  "
Set-TextValue $ws "D27" "71.31"
$ws.Range("E27").Value = "  +1.00%  "
Set-TextValue $ws "D28" "1.68"
$ws.Range("E28").Value = "  -0.44%  "
Set-TextValue $ws "D29" "7.58"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws "D31" "8.27"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D32" "0.161"
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("E34").Value = "  +0.04%  "
Set-TextValue $ws "D35" "23.47"
$ws.Range("E35").Value = "  +0.44%  "
Set-TextValue $ws "D36" "5.36"
$ws.Range("E36").Value = "  -3.80%  "
Set-TextValue $ws "D37" "1.56"
$ws.Range("E37").Value = "  +1.19%  "
Set-TextValue $ws "D38" "6.87"
$ws.Range("E38").Value = "  -0.78%  "
Set-TextValue $ws "D39" "165.35"
$ws.Range("E39").Value = "  +1.36%  "
Set-TextValue $ws "D40" "0.0784"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D41" "1.74"
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("B42").Value = "ONDO"
$ws.Range("C42").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue $ws "D42" "1.24"
$ws.Range("E42").Value = "  +2.73%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D43" "0.783"
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D44" "0.998"
$ws.Range("E44").Value = "  -0.19%  "
Set-TextValue $ws "D45" "25.25"
$ws.Range("E45").Value = "  +9.66%  "
Set-TextValue $ws "D46" "4.43"
$ws.Range("E46").Value = "  +0.69%  "
Set-TextValue $ws "D47" "41.53"
$ws.Range("E47").Value = "  +0.51%  "
Set-TextValue $ws "D48" "6.87"
$ws.Range("E48").Value = "  -1.52%  "
Set-TextValue $ws "D49" "22.84"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "2.342.00"
$ws.Range("E50").Value = "  +6.02%  "
Set-TextValue $ws "D51" "0.0262"
$ws.Range("E51").Value = "  +2.38%  "
